# Add "return from movie page" button -> in data terms: fill in the missing
# dis_link/img columns (E/F) for several Fantasyland attractions, and fix the
# land label for the Tomorrowland attractions that were mistakenly tagged as
# "Fantasyland".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Introduce the new "Tomorrowland" land label first (so it becomes the
#    first newly-appended shared string), fixing rows 29-35 which were all
#    mislabeled as "Fantasyland".
$ws.Range("C29").Value = "Tomorrowland"
$ws.Range("C30").Value = "Tomorrowland"
$ws.Range("C31").Value = "Tomorrowland"
$ws.Range("C32").Value = "Tomorrowland"
$ws.Range("C33").Value = "Tomorrowland"
$ws.Range("C34").Value = "Tomorrowland"
$ws.Range("C35").Value = "Tomorrowland"

# 2) Fill in dis_link (E) / img (F) columns for the Fantasyland attractions
#    that didn't have them yet.

# Row 16 - Dumbo the Flying Elephant
$ws.Range("E16").Value = "https://disneyworld.disney.go.com/attractions/magic-kingdom/dumbo-the-flying-elephant/"
$ws.Range("F16").Value = "https://cdn1.parksmedia.wdprapps.disney.com/resize/mwImage/1/1600/900/75/dam/wdpro-assets/gallery/attractions/magic-kingdom/dumbo-the-flying-elephant/dumbo-the-flying-elephant-gallery06.jpg?1550820961335"

# Row 17 - Enchanted Tales with Belle
$ws.Range("E17").Value = "https://disneyworld.disney.go.com/attractions/magic-kingdom/enchanted-tales-with-belle/"
$ws.Range("F17").Value = "https://www.tripsavvy.com/thmb/TRwzYMvBr04zEAap2z23uDDPVjM=/2700x1887/filters:fill(auto,1)/Enchanted-Tales-with-Belle-Disney-World-58bdeef73df78c353cdda5cd.jpg"

# Row 18 - It's a Small World (image string introduced before the link string)
$ws.Range("F18").Value = "https://cdn1.parksmedia.wdprapps.disney.com/resize/mwImage/1/1600/900/75/dam/wdpro-assets/parks-and-tickets/attractions/magic-kingdom/its-a-small-world/its-a-small-world-00.jpg?1634639538944"
$ws.Range("E18").Value = "https://disneyworld.disney.go.com/attractions/magic-kingdom/its-a-small-world/"

# Row 19 - Mad Tea Party (image string introduced before the link string)
$ws.Range("F19").Value = "https://cdn1.parksmedia.wdprapps.disney.com/resize/mwImage/1/1600/900/75/dam/wdpro-assets/parks-and-tickets/attractions/magic-kingdom/mad-tea-party/mad-tea-party-00.jpg?1634619687993"
$ws.Range("E19").Value = "https://disneyworld.disney.go.com/attractions/magic-kingdom/mad-tea-party/"

# Row 20 - Mickey's PhilharMagic
$ws.Range("E20").Value = "https://disneyworld.disney.go.com/attractions/magic-kingdom/mickeys-philharmagic/"
$ws.Range("F20").Value = "https://cdn1.parksmedia.wdprapps.disney.com/resize/mwImage/1/1600/900/75/dam/wdpro-assets/gallery/attractions/magic-kingdom/mickeys-philharmagic/mickeys-philharmagic-gallery01.jpg?1551208516468"

# Row 21 - Peter Pan's Flight (image string introduced before the link string)
$ws.Range("F21").Value = "https://cdn1.parksmedia.wdprapps.disney.com/resize/mwImage/1/1600/900/75/dam/wdpro-assets/gallery/attractions/magic-kingdom/peter-pans-flight/peter-pans-flight-gallery03.jpg?1559892376214"
$ws.Range("E21").Value = "https://disneyworld.disney.go.com/attractions/magic-kingdom/peter-pan-flight/"

# 3) Leave the cursor/selection where the author left it.
$ws.Range("D22").Select()
